$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 0.1884
$ws.Range("C5").Value = 0.0215
$ws.Range("D5").Value = 0.131
$ws.Range("E5").Value = 0.2764
$ws.Range("H5").Value = -0.0847
$ws.Range("L5").Value = -0.1382
$ws.Range("C6").Value = -0.024
$ws.Range("D6").Value = -0.0296
$ws.Range("E6").Value = -0.025
$ws.Range("F6").Value = -0.0515
$ws.Range("G6").Value = -0.0491
$ws.Range("H6").Value = -0.055
$ws.Range("I6").Value = -0.0401
$ws.Range("J6").Value = -0.0723
$ws.Range("K6").Value = -0.068
$ws.Range("L6").Value = -0.0427
$ws.Range("M6").Value = -0.0346
$ws.Range("C7").Value = -0.2563
$ws.Range("D7").Value = -0.0177
$ws.Range("E7").Value = -0.0188
$ws.Range("F7").Value = 0.1253
$ws.Range("G7").Value = 0.2236
$ws.Range("H7").Value = -0.0851
$ws.Range("I7").Value = -0.1174
$ws.Range("J7").Value = -0.1602
$ws.Range("K7").Value = -0.182
$ws.Range("L7").Value = 0.043
$ws.Range("M7").Value = 0.1276
$ws.Range("C8").Value = -0.2721
$ws.Range("D8").Value = -0.196
$ws.Range("E8").Value = -0.1689
$ws.Range("F8").Value = -0.1238
$ws.Range("G8").Value = -0.3971
$ws.Range("H8").Value = -0.303
$ws.Range("I8").Value = -0.2026
$ws.Range("J8").Value = -0.0789
$ws.Range("K8").Value = -0.0342
$ws.Range("L8").Value = -0.0175
$ws.Range("M8").Value = 0.002
$ws.Range("D9").Value = 0.5001
$ws.Range("D10").Value = 0.4835
$ws.Range("C11").Value = 0.0102
$ws.Range("C12").Value = -0.5948
$ws.Range("D12").Value = -1.0492
$ws.Range("G12").Value = -0.5319
$ws.Range("C13").Value = 0.3833
$ws.Range("D13").Value = 0.2726
$ws.Range("E13").Value = -0.1238
$ws.Range("F13").Value = -0.0875
$ws.Range("G13").Value = -0.1085
$ws.Range("H13").Value = -0.0991
$ws.Range("I13").Value = -0.0592
$ws.Range("J13").Value = -0.048
$ws.Range("L13").Value = -0.6506
$ws.Range("M13").Value = -0.3081
$ws.Range("C14").Value = 0.1724
$ws.Range("D14").Value = -0.7572
$ws.Range("E14").Value = -1.1114
$ws.Range("F14").Value = -0.9955
$ws.Range("G14").Value = -1.4852
$ws.Range("H14").Value = -1.0685
$ws.Range("I14").Value = -0.4401
$ws.Range("J14").Value = -0.4047
$ws.Range("K14").Value = -0.321
$ws.Range("L14").Value = -0.1744
$ws.Range("C15").Value = -2.7628
$ws.Range("D15").Value = -1.9883
$ws.Range("E15").Value = -2.9886
$ws.Range("F15").Value = -2.1825
$ws.Range("G15").Value = -3.7821
$ws.Range("H15").Value = -2.0184
$ws.Range("I15").Value = -1.6445
$ws.Range("J15").Value = -1.7017
$ws.Range("K15").Value = -2.8308
$ws.Range("L15").Value = -1.1846
$ws.Range("M15").Value = -0.5898
$ws.Range("C16").Value = 22741
$ws.Range("C17").Value = 0.5458
$ws.Range("D17").Value = 0.2185
$ws.Range("E17").Value = 0.3464
$ws.Range("H17").Value = 0.1932
$ws.Range("I17").Value = 0.2091
$ws.Range("K17").Value = -0.2326
$ws.Range("C19").Value = -1.0494
$ws.Range("D19").Value = 0.0464
$ws.Range("E19").Value = -1.2305
$ws.Range("F19").Value = 0.0112
$ws.Range("G19").Value = -0.0313
$ws.Range("H19").Value = -0.0333
$ws.Range("I19").Value = -0.0353
$ws.Range("J19").Value = -0.0348
$ws.Range("K19").Value = -1.0901
$ws.Range("L19").Value = -0.0005
$ws.Range("M19").Value = -0.0307
$ws.Range("C20").Value = -0.2711
$ws.Range("D20").Value = -1.2688
$ws.Range("E20").Value = -0.2747
$ws.Range("G20").Value = -1.0259
$ws.Range("C21").Value = -0.7397
$ws.Range("D21").Value = 0.1726
$ws.Range("E21").Value = 0.1617
$ws.Range("K21").Value = 0.1831
$ws.Range("M21").Value = -0.1338
$ws.Range("C22").Value = -0.014
$ws.Range("J22").Value = -0.0222
$ws.Range("C23").Value = 0.1964
$ws.Range("D23").Value = -0.0326
$ws.Range("E23").Value = 0.1275
$ws.Range("F23").Value = 0.0049
$ws.Range("G23").Value = -0.0708
$ws.Range("H23").Value = 0.1439
$ws.Range("I23").Value = 0.1399
$ws.Range("J23").Value = 0.2038
$ws.Range("K23").Value = 0.1915
$ws.Range("L23").Value = -0.0178
$ws.Range("M23").Value = 0.0084
$ws.Range("C24").Value = -0.0802
$ws.Range("D24").Value = -0.1476
$ws.Range("E24").Value = -0.1458
$ws.Range("F24").Value = -0.1149
$ws.Range("G24").Value = -0.1726
$ws.Range("H24").Value = -0.128
$ws.Range("I24").Value = -0.1333
$ws.Range("J24").Value = -0.1044
$ws.Range("K24").Value = -0.0698
$ws.Range("L24").Value = -0.0311
$ws.Range("M24").Value = -0.0188
$ws.Range("C25").Value = 0.0178
$ws.Range("C28").Value = -0.2286
$ws.Range("D28").Value = -0.2414
$ws.Range("E28").Value = -0.0821
$ws.Range("F28").Value = -0.119
$ws.Range("G28").Value = -0.2302
$ws.Range("H28").Value = -0.077
$ws.Range("I28").Value = -0.0404
$ws.Range("J28").Value = -0.0391
$ws.Range("K28").Value = -0.0247
$ws.Range("L28").Value = -0.0089
$ws.Range("M28").Value = -0.0042
$ws.Range("D31").Value = -0.0001
$ws.Range("C32").Value = 0.0066
$ws.Range("D32").Value = 0.1623
$ws.Range("E32").Value = -0.1753
$ws.Range("C33").Value = -0.0084
$ws.Range("D33").Value = -0.0085
$ws.Range("E33").Value = -0.0083
$ws.Range("F33").Value = 0.0002
$ws.Range("G33").Value = 0.0002
$ws.Range("H33").Value = 0.0002
$ws.Range("I33").Value = 0.0002
$ws.Range("J33").Value = 0.0002
$ws.Range("K33").Value = 0.0002
$ws.Range("L33").Value = 0.0002
$ws.Range("M33").Value = 0.0002
$ws.Range("C34").Value = 0.0001
$ws.Range("D34").Value = -0.0338
$ws.Range("E34").Value = -0.0323
$ws.Range("F34").Value = -0.0465
$ws.Range("G34").Value = -0.0452
$ws.Range("H34").Value = -0.0094
$ws.Range("I34").Value = -0.0098
$ws.Range("J34").Value = 0.0045
$ws.Range("K34").Value = 0.0039
$ws.Range("L34").Value = 0.0004
$ws.Range("M34").Value = -0.0007
$ws.Range("C35").Value = -0.0033
$ws.Range("D35").Value = -0.0031
$ws.Range("E35").Value = -0.0015
$ws.Range("F35").Value = -0.0011
$ws.Range("G35").Value = -0.0011
$ws.Range("H35").Value = -0.0013
$ws.Range("I35").Value = -0.0013
$ws.Range("J35").Value = -0.0013
$ws.Range("K35").Value = 0.0002
$ws.Range("L35").Value = 0.0002
$ws.Range("M35").Value = 0.0002
$ws.Range("D36").Value = -0.0002
$ws.Range("D37").Value = -0.0002
$ws.Range("C38").Value = 0.0003
$ws.Range("C39").Value = 0.021
$ws.Range("D39").Value = 0.0004
$ws.Range("G39").Value = -0.0192
$ws.Range("C40").Value = 0.012
$ws.Range("D40").Value = -0.0789
$ws.Range("E40").Value = 0.0103
$ws.Range("F40").Value = 0.0034
$ws.Range("G40").Value = 0.0033
$ws.Range("H40").Value = 0.0037
$ws.Range("I40").Value = 0.0004
$ws.Range("J40").Value = 0.0004
$ws.Range("K40").Value = -0.0001
$ws.Range("L40").Value = 0.0035
$ws.Range("M40").Value = 0.0004
$ws.Range("D41").Value = 0.148
$ws.Range("E41").Value = -0.0229
$ws.Range("F41").Value = -0.0082
$ws.Range("G41").Value = -0.008
$ws.Range("H41").Value = -0.042
$ws.Range("I41").Value = -0.0276
$ws.Range("J41").Value = -0.0003
$ws.Range("K41").Value = -0.007
$ws.Range("L41").Value = -0.0266
$ws.Range("C42").Value = 0.0653
$ws.Range("D42").Value = 0.0855
$ws.Range("E42").Value = -0.0846
$ws.Range("F42").Value = -0.0452
$ws.Range("G42").Value = -0.1027
$ws.Range("H42").Value = 0.0094
$ws.Range("I42").Value = 0.0117
$ws.Range("J42").Value = 0.0518
$ws.Range("K42").Value = -0.0097
$ws.Range("L42").Value = -0.0183
$ws.Range("M42").Value = -0.0128
$ws.Range("C43").Value = 9.6
$ws.Range("C44").Value = 0.0161
$ws.Range("D44").Value = 0.1623
$ws.Range("E44").Value = -0.1753
$ws.Range("C46").Value = 0.0001
$ws.Range("D46").Value = 0.0461
$ws.Range("E46").Value = -0.0132
$ws.Range("F46").Value = -0.0003
$ws.Range("G46").Value = -0.062
$ws.Range("H46").Value = 0.0134
$ws.Range("I46").Value = 0.0002
$ws.Range("J46").Value = 0.0162
$ws.Range("K46").Value = 0.0001
$ws.Range("M46").Value = -0.0145
$ws.Range("C47").Value = 0.0001
$ws.Range("D47").Value = 0.0006
$ws.Range("C48").Value = -0.0307
$ws.Range("D48").Value = -0.1625
$ws.Range("E48").Value = 0.1753
$ws.Range("C49").Value = -0.0008
$ws.Range("C50").Value = -0.0204
$ws.Range("D50").Value = -0.0285
$ws.Range("E50").Value = -0.0276
$ws.Range("F50").Value = -0.0114
$ws.Range("G50").Value = 0.0076
$ws.Range("H50").Value = 0.0128
$ws.Range("I50").Value = 0.0126
$ws.Range("J50").Value = -0.0023
$ws.Range("K50").Value = -0.0024
$ws.Range("L50").Value = 0.0004
$ws.Range("M50").Value = -0.0001
$ws.Range("C51").Value = 0.0892
$ws.Range("D51").Value = 0.0803
$ws.Range("E51").Value = 0.0362
$ws.Range("F51").Value = 0.0389
$ws.Range("G51").Value = 0.0384
$ws.Range("H51").Value = 0.0379
$ws.Range("I51").Value = 0.0374
$ws.Range("J51").Value = 0.037
$ws.Range("K51").Value = -0.0023
$ws.Range("L51").Value = 0.0005
$ws.Range("M51").Value = 0.0003
$ws.Range("C52").Value = 0.0003
$ws.Range("C55").Value = -0.0004
$ws.Range("D55").Value = -0.0365
$ws.Range("E55").Value = -0.0253
$ws.Range("F55").Value = -0.0202
$ws.Range("G55").Value = -0.0059
$ws.Range("H55").Value = -0.0059
$ws.Range("I55").Value = -0.0005
$ws.Range("J55").Value = -0.0026
$ws.Range("K55").Value = -0.0024
$ws.Range("L55").Value = 0.003
$ws.Range("M55").Value = 0.0014
